# 6-diff-tdd-bdd.pptx - "docs: new ver on pdfs and folders"
#
# 1) Bump the cached "datetimeFigureOut" field text from 2024-10-26 to
#    2024-10-27 everywhere it appears: the slide master's Date
#    Placeholder and the Date Placeholder on every one of the 11
#    slide layouts.
# 2) Remove the extra "TextBox 4" shape (free-floating duplicate of the
#    title, "Syfte och Fokus") that was left on slide 2.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "2024-10-26") {
                $tr.Text = "2024-10-27"
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 2: delete the stray "TextBox 4" ("Syfte och Fokus") shape.
$slide2 = $p.Slides.Item(2)
for ($i = $slide2.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 4") {
        $shape.Delete()
    }
}
